# The deck's slide design is switched from the custom "Integral" theme to
# the default Office Theme color palette (ppt/theme/theme1.xml, which is
# the theme referenced by the slide master).
#
# PowerPoint COM exposes the 12 theme colors (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) through SlideMaster.Theme.ThemeColorScheme, indexed 1-12
# in that order. Each item's .RGB is a plain Win32 COLORREF
# (0x00BBGGRR === R + G*256 + B*65536), so we convert the target hex swatches
# accordingly.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

function HexToCOMRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme color scheme (the new design), in MsoThemeColorSchemeIndex order.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $tcs.Item($i).RGB = HexToCOMRGB $officeColors[$i - 1]
}
